{"js": "// Replace the three-digit x one-digit multiplication expressions in the\n// document's table cells with their updated values, per the commit diff.\n// Each \"NNN\u00d7N=\" string is unique in the document, so a plain text search\n// and replace (preserving the existing run formatting) is sufficient.\n\nconst replacements = [\n  [\"478\u00d74=\", \"497\u00d72=\"],\n  [\"892\u00d77=\", \"254\u00d74=\"],\n  [\"633\u00d78=\", \"284\u00d76=\"],\n  [\"791\u00d77=\", \"245\u00d77=\"],\n  [\"389\u00d78=\", \"994\u00d75=\"],\n  [\"619\u00d77=\", \"578\u00d73=\"],\n  [\"222\u00d73=\", \"754\u00d78=\"],\n  [\"992\u00d74=\", \"338\u00d79=\"],\n  [\"613\u00d74=\", \"441\u00d78=\"],\n  [\"739\u00d75=\", \"742\u00d75=\"],\n  [\"190\u00d72=\", \"684\u00d78=\"],\n  [\"517\u00d76=\", \"193\u00d79=\"],\n  [\"752\u00d76=\", \"435\u00d79=\"],\n  [\"284\u00d72=\", \"609\u00d78=\"],\n  [\"488\u00d77=\", \"395\u00d77=\"],\n  [\"437\u00d78=\", \"646\u00d75=\"],\n  [\"788\u00d73=\", \"106\u00d72=\"],\n  [\"227\u00d75=\", \"636\u00d72=\"],\n  [\"564\u00d77=\", \"862\u00d76=\"],\n  [\"376\u00d73=\", \"911\u00d72=\"],\n  [\"114\u00d78=\", \"686\u00d73=\"],\n  [\"255\u00d78=\", \"933\u00d74=\"],\n  [\"407\u00d76=\", \"549\u00d73=\"],\n  [\"816\u00d76=\", \"817\u00d74=\"],\n  [\"841\u00d73=\", \"345\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication expressions in the\n# document's table cells with their updated values, per the commit diff.\n# Each \"NNN\u00d7N=\" string is unique in the document, so a plain Find/Replace\n# (preserving the existing run formatting) is sufficient.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"478\u00d74=\", \"497\u00d72=\"),\n    @(\"892\u00d77=\", \"254\u00d74=\"),\n    @(\"633\u00d78=\", \"284\u00d76=\"),\n    @(\"791\u00d77=\", \"245\u00d77=\"),\n    @(\"389\u00d78=\", \"994\u00d75=\"),\n    @(\"619\u00d77=\", \"578\u00d73=\"),\n    @(\"222\u00d73=\", \"754\u00d78=\"),\n    @(\"992\u00d74=\", \"338\u00d79=\"),\n    @(\"613\u00d74=\", \"441\u00d78=\"),\n    @(\"739\u00d75=\", \"742\u00d75=\"),\n    @(\"190\u00d72=\", \"684\u00d78=\"),\n    @(\"517\u00d76=\", \"193\u00d79=\"),\n    @(\"752\u00d76=\", \"435\u00d79=\"),\n    @(\"284\u00d72=\", \"609\u00d78=\"),\n    @(\"488\u00d77=\", \"395\u00d77=\"),\n    @(\"437\u00d78=\", \"646\u00d75=\"),\n    @(\"788\u00d73=\", \"106\u00d72=\"),\n    @(\"227\u00d75=\", \"636\u00d72=\"),\n    @(\"564\u00d77=\", \"862\u00d76=\"),\n    @(\"376\u00d73=\", \"911\u00d72=\"),\n    @(\"114\u00d78=\", \"686\u00d73=\"),\n    @(\"255\u00d78=\", \"933\u00d74=\"),\n    @(\"407\u00d76=\", \"549\u00d73=\"),\n    @(\"816\u00d76=\", \"817\u00d74=\"),\n    @(\"841\u00d73=\", \"345\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
